$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and the SEI/THORChain row swap)
# Force text format on the cells we touch so numeric-looking strings (e.g. "3.30",
# "0.0830") keep their exact original text representation instead of being
# reinterpreted/normalized as numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.817.66'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.791.34'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.09%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '354.29'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.30%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.34'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.11%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.551'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.52%  '

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.57%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.99'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.41%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.01%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.24'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.67%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0839'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.93%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.64'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.31%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.233.00'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.16%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.815.66'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.17%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.929'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.69%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.793.05'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.22%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.44%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.12'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.34%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.17'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.68%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0966'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.54%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.06'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.29%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '266.56'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.95%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.73'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.56%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.18'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.97%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.01%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.162'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +11.75%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.24'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.46%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.97'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +8.22%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.23'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +8.77%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '51.99'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.82%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.89%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.60'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.20%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -12.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0830'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.83%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.57'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.85%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.16'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.32%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.97'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.38%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.56'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.40%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.78%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '120.94'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.06'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.43%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.55%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.137.42'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.01%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.30'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.27%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.33'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +7.46%  '

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.45'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.89%  '

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'SEI'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.911'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.19%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +9.30%  '
